$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws.Range("B2").Value = 0.3161506809260008
$ws.Range("B3").Value = 63.24789231130328
$ws.Range("B4").Value = 63.52900847697614
$ws.Range("B5").Value = 63.49399691014947
$ws.Range("B6").Value = 64.48861907087098
$ws.Range("B7").Value = 65.22023036518915
$ws.Range("B8").Value = 66.12891524107229
$ws.Range("B9").Value = 66.33679213525222
$ws.Range("B10").Value = 67.38243195902817
$ws.Range("B11").Value = 68.64901148156481
$ws.Range("B12").Value = 69.03150881498213
$ws.Range("B13").Value = 69.30583316918772
$ws.Range("B14").Value = 69.79011538731207
$ws.Range("B15").Value = 69.76594441009851
$ws.Range("B16").Value = 69.87656780615016
$ws.Range("B17").Value = 70.57760914490584
$ws.Range("B18").Value = 72.22685714459905
$ws.Range("B19").Value = 70.86034721423181
$ws.Range("B20").Value = 72.00727312586457
$ws.Range("B21").Value = 72.06327396208286
$ws.Range("B22").Value = 71.93692329255649
$ws.Range("B23").Value = 71.20017202974088
$ws.Range("B24").Value = 70.11168154403589

$ws = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws.Range("B2").Value = 0.2994451769062907
$ws.Range("B3").Value = 63.23118680728358
$ws.Range("B4").Value = 63.51230297295643
$ws.Range("B5").Value = 63.47729140612976
$ws.Range("B6").Value = 64.47191356685127
$ws.Range("B7").Value = 65.20352486116944
$ws.Range("B8").Value = 66.11220973705258
$ws.Range("B9").Value = 66.32008663123251
$ws.Range("B10").Value = 67.36572645500846
$ws.Range("B11").Value = 68.6323059775451
$ws.Range("B12").Value = 69.01480331096242
$ws.Range("B13").Value = 69.28912766516801
$ws.Range("B14").Value = 69.77340988329236
$ws.Range("B15").Value = 69.7492389060788
$ws.Range("B16").Value = 69.85986230213045
$ws.Range("B17").Value = 70.56090364088612
$ws.Range("B18").Value = 72.21015164057934
$ws.Range("B19").Value = 70.8436417102121
$ws.Range("B20").Value = 71.99056762184486
$ws.Range("B21").Value = 72.04656845806315
$ws.Range("B22").Value = 71.92021778853677
$ws.Range("B23").Value = 71.18346652572117
$ws.Range("B24").Value = 70.11168154403589
$ws.Range("B25").Value = 70.18650249106669
$ws.Range("B26").Value = 69.83522074360363
$ws.Range("B27").Value = 69.73196446797155
$ws.Range("B28").Value = 69.71618800114416
$ws.Range("B29").Value = 69.76484494805047
